# Update countries & provincias Spain
# Applies the 23-Oct-2020 data refresh (08:14 -> 09:31) to the "Pais" sheet:
#  - refreshed case counters for several countries
#  - re-sorted rows that swapped rank because of the refreshed counters
#    (Armenia/Austria, Hungria/Libia, Georgia/Malasia, Montserrat/Islas Malvinas)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp caption (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 09:31"

# --- Estados Unidos (row 4) : values refreshed, no re-sort needed ---
$ws.Range("B4").Value = 8661722
$ws.Range("C4").Value = 71
$ws.Range("D4").Value = 5655327
$ws.Range("E4").Value = 2778014

# --- Armenia overtakes Austria -> rows 59/60 swap ---
$ws.Range("A59").Value = "Armenia"
$ws.Range("B59").Value = 73310
$ws.Range("C59").Value = 2474
$ws.Range("D59").Value = 50276
$ws.Range("E59").Value = 21889
$ws.Range("G59").Value = 14
$ws.Range("H59").Value = 1145

$ws.Range("A60").Value = "Austria"
$ws.Range("B60").Value = 71844
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 53970
$ws.Range("E60").Value = 16933
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 941

# --- Hungria overtakes Libia -> rows 70/71 swap ---
$ws.Range("A70").Value = "Hungria"
$ws.Range("B70").Value = 54278
$ws.Range("C70").Value = 2066
$ws.Range("D70").Value = 15655
$ws.Range("E70").Value = 37271
$ws.Range("G70").Value = 47
$ws.Range("H70").Value = 1352

$ws.Range("A71").Value = "Libia"
$ws.Range("B71").Value = 52620
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 29057
$ws.Range("E71").Value = 22795
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 768

# --- Afganistan (row 79) : values refreshed, no re-sort needed ---
$ws.Range("B79").Value = 40687
$ws.Range("C79").Value = 61
$ws.Range("D79").Value = 34010
$ws.Range("E79").Value = 5170
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 1507

# --- Georgia overtakes Malasia -> rows 91/92 swap ---
$ws.Range("A91").Value = "Georgia"
$ws.Range("B91").Value = 24562
$ws.Range("C91").Value = 1759
$ws.Range("D91").Value = 9751
$ws.Range("E91").Value = 14628
$ws.Range("G91").Value = 5
$ws.Range("H91").Value = 183

$ws.Range("A92").Value = "Malasia"
$ws.Range("B92").Value = 23804
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 15417
$ws.Range("E92").Value = 8183
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 204

# --- Montserrat overtakes Islas Malvinas -> rows 216/217 swap ---
$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1

$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 13
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 0
